$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agendas")

# Clear the old "s" / "09:00-13:00" entry that lived under the Lunes column
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Store the new available advisory hour block under the Viernes column
$ws.Range("F2").Value = "08:00 - 09:00"
